# Applies the "Upload new version with timestamp" edit:
#  - Inserts a new item row for "DOLIPRANE 1 GM 15 TABS." right before the
#    "IVYMOND SYRUP" row (new row 11), shifting all following item rows,
#    the totals row and the footer row down by one.
#  - Renumbers the serial (م) column for every item row pushed down.
#  - Refreshes the totals (سعر البيع sum) cell for the now-longer list.
#  - Updates the generated timestamp shown in the footer.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a blank row at row 11 (everything from row 11 down shifts to +1).
$ws.Rows("11:11").Insert()

# 2. Use the row that is now at 12 (the old "IVYMOND SYRUP" row) as a
#    formatting/merge template for the freshly inserted blank row 11, then
#    overwrite its contents with the new product's data.
$ws.Range("A12:Q12").Copy($ws.Range("A11:Q11"))
$excel.CutCopyMode = $false

$ws.Range("A11").Value = 5
$ws.Range("C11").Value = "DOLIPRANE 1 GM 15 TABS."
$ws.Range("H11").Value = "'6:0"
$ws.Range("L11").Value = "'1"
$ws.Range("N11").Value = "'48.00"
$ws.Range("P11").Value = "'15.8400"
$ws.Range("Q11").Value = "'0:1"

# 3. Renumber the serial column (م) for every item row that got pushed
#    down by the insert (old rows 11-19 are now rows 12-20).
$ws.Range("A12").Value = 6
$ws.Range("A13").Value = 7
$ws.Range("A14").Value = 8
$ws.Range("A15").Value = 9
$ws.Range("A16").Value = 10
$ws.Range("A17").Value = 11
$ws.Range("A18").Value = 12
$ws.Range("A19").Value = 13
$ws.Range("A20").Value = 14

# 4. Refresh the "سعر البيع" total now that a new row participates in it.
$ws.Range("P21").Value = 671.75

# 5. Bump the generated timestamp shown in the footer.
$ws.Range("A22").Value = "Wednesday, 10 September, 2025 12:43 PM"
